$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before column A, shifting the original table right ---
$ws.Range("A:B").Insert()
$ws.Range("B1").Value = 0

# --- Row 4: header row (# droplet) ---
$ws.Range("A4").Value = "# droplet"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 5
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 7
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 9
$ws.Range("L4").Value = 10
$ws.Range("N4").Value = "Average"

# --- Row 5: merged sub-header label (1st part test) ---
$ws.Range("B5").Value = "Ave Time of the Detection Part (pre + detect + post)"

# --- Row 6: Run Time (ms) values for the 1st part test ---
$ws.Range("A6").Value = "Run Time (ms)"
$ws.Range("B6").Value = 353.70377999999988
$ws.Range("C6").Value = 357.05817333333317
$ws.Range("D6").Value = 343.73861999999986
$ws.Range("E6").Value = 344.2699933333331
$ws.Range("F6").Value = 347.89613999999972
$ws.Range("G6").Value = 351.1491666666663
$ws.Range("H6").Value = 355.95284666666663
$ws.Range("I6").Value = 344.48259333333311
$ws.Range("J6").Value = 339.54561333333305
$ws.Range("K6").Value = 344.88387333333333
$ws.Range("L6").Value = 348.26161999999977
$ws.Range("N6").Formula = "=AVERAGE(B6:L6)"

# --- Row 7: merged sub-header label (2nd part test) ---
$ws.Range("B7").Value = "Average Time of the YOLOv5 Detection (only detect)"

# --- Row 8: Run Time (ms) formulas for the 2nd part test ---
$ws.Range("A8").Value = "Run Time (ms)"
$ws.Range("B8").Formula = "=1.6+67.5+1"
$ws.Range("C8").Formula = "=0.9+67.7+0.7"
$ws.Range("D8").Formula = "=0.9+67.4+0.8"
$ws.Range("E8").Formula = "=1+66.6+1.1"
$ws.Range("F8").Formula = "=0.6+68.4+0.8"
$ws.Range("G8").Formula = "=0.7+69.4+0.5"
$ws.Range("H8").Formula = "=0.9+69.8+0.5"
$ws.Range("I8").Formula = "=0.7+66+0.8"
$ws.Range("K8").Formula = "=1+69.1+0.7"
$ws.Range("L8").Formula = "=0.7+70+0.9"
$ws.Range("N8").Formula = "=(SUM(B8:I8) + SUM(K8:L8)) / 10"

# --- Alignment: centre the data rows (including column A on every row 4-8) ---
$ws.Range("A4:N4").HorizontalAlignment = -4108
$ws.Range("A4:N4").VerticalAlignment = -4108
$ws.Range("A5").HorizontalAlignment = -4108
$ws.Range("A5").VerticalAlignment = -4108
$ws.Range("A6:N6").HorizontalAlignment = -4108
$ws.Range("A6:N6").VerticalAlignment = -4108
$ws.Range("A7").HorizontalAlignment = -4108
$ws.Range("A7").VerticalAlignment = -4108
$ws.Range("A8:N8").HorizontalAlignment = -4108
$ws.Range("A8:N8").VerticalAlignment = -4108

# --- Merge and centre the sub-header label rows ---
$ws.Range("B5:N5").Merge()
$ws.Range("B7:N7").Merge()
$ws.Range("B5:N5").HorizontalAlignment = -4108
$ws.Range("B5:N5").VerticalAlignment = -4108
$ws.Range("B7:N7").HorizontalAlignment = -4108
$ws.Range("B7:N7").VerticalAlignment = -4108

# --- Column A width (auto-fit to the longest label, "# droplet") ---
$ws.Columns("A:A").AutoFit()

$ws.Range("O17").Select()
